{"js": "// Update the answer cells of the \"two-digit \u00f7 one-digit\" practice table.\n// Each non-empty cell in the (single) table holds one division problem in\n// document order; replace the text of each one with its new value while\n// leaving every other part of the document (formatting, empty rows, the\n// date paragraph, ...) untouched.\n\nconst oldNewPairs = [\n  [\"99\u00f78=12, 3\", \"62\u00f72=31, 0\"],\n  [\"82\u00f75=16, 2\", \"89\u00f77=12, 5\"],\n  [\"29\u00f75=5, 4\", \"53\u00f74=13, 1\"],\n  [\"71\u00f77=10, 1\", \"67\u00f73=22, 1\"],\n  [\"78\u00f73=26, 0\", \"68\u00f79=7, 5\"],\n  [\"21\u00f75=4, 1\", \"31\u00f74=7, 3\"],\n  [\"61\u00f74=15, 1\", \"64\u00f76=10, 4\"],\n  [\"54\u00f77=7, 5\", \"74\u00f74=18, 2\"],\n  [\"40\u00f72=20, 0\", \"53\u00f78=6, 5\"],\n  [\"79\u00f73=26, 1\", \"59\u00f79=6, 5\"],\n  [\"93\u00f77=13, 2\", \"22\u00f76=3, 4\"],\n  [\"78\u00f73=26, 0\", \"57\u00f77=8, 1\"],\n  [\"97\u00f73=32, 1\", \"46\u00f75=9, 1\"],\n  [\"69\u00f73=23, 0\", \"78\u00f72=39, 0\"],\n  [\"88\u00f73=29, 1\", \"20\u00f79=2, 2\"],\n  [\"17\u00f73=5, 2\", \"82\u00f79=9, 1\"],\n  [\"52\u00f74=13, 0\", \"58\u00f72=29, 0\"],\n  [\"49\u00f72=24, 1\", \"28\u00f78=3, 4\"],\n  [\"31\u00f76=5, 1\", \"29\u00f72=14, 1\"],\n  [\"61\u00f73=20, 1\", \"86\u00f73=28, 2\"],\n  [\"29\u00f75=5, 4\", \"77\u00f74=19, 1\"],\n  [\"89\u00f79=9, 8\", \"13\u00f79=1, 4\"],\n  [\"89\u00f75=17, 4\", \"71\u00f76=11, 5\"],\n  [\"93\u00f77=13, 2\", \"68\u00f75=13, 3\"],\n  [\"98\u00f72=49, 0\", \"42\u00f76=7, 0\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table with the answer cells, found none.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nlet pairIndex = 0;\nconst updatedValues = table.values.map((row) =>\n  row.map((cellText) => {\n    if (cellText === \"\" || cellText === null || cellText === undefined) {\n      // Blank spacer cell/row - leave as-is.\n      return cellText;\n    }\n    if (pairIndex >= oldNewPairs.length) {\n      return cellText;\n    }\n    const [oldValue, newValue] = oldNewPairs[pairIndex];\n    pairIndex++;\n    if (cellText !== oldValue) {\n      throw new Error(\n        \"Unexpected cell text at position \" +\n          pairIndex +\n          \": expected '\" +\n          oldValue +\n          \"' but found '\" +\n          cellText +\n          \"'.\"\n      );\n    }\n    return newValue;\n  })\n);\n\nif (pairIndex !== oldNewPairs.length) {\n  throw new Error(\n    \"Only matched \" + pairIndex + \" of \" + oldNewPairs.length + \" expected cells.\"\n  );\n}\n\ntable.values = updatedValues;\nawait context.sync();\n", "ps1": "# Update each division-problem answer cell in the practice table.\n# The divide sign (U+00F7) is built from its code point rather than\n# embedded literally, and every \"+\" concatenation is forced to string\n# concatenation via an explicit [string] cast (otherwise a leading\n# numeric-looking operand gets added to the character code instead of\n# concatenated with it).\n$div = [string][char]0xF7\n\n$entries = @(\n  @{ row=1; col=1; leftOld=\"99\"; rightOld=\"8=12, 3\"; leftNew=\"62\"; rightNew=\"2=31, 0\" },\n  @{ row=1; col=2; leftOld=\"82\"; rightOld=\"5=16, 2\"; leftNew=\"89\"; rightNew=\"7=12, 5\" },\n  @{ row=1; col=3; leftOld=\"29\"; rightOld=\"5=5, 4\"; leftNew=\"53\"; rightNew=\"4=13, 1\" },\n  @{ row=1; col=4; leftOld=\"71\"; rightOld=\"7=10, 1\"; leftNew=\"67\"; rightNew=\"3=22, 1\" },\n  @{ row=1; col=5; leftOld=\"78\"; rightOld=\"3=26, 0\"; leftNew=\"68\"; rightNew=\"9=7, 5\" },\n  @{ row=5; col=1; leftOld=\"21\"; rightOld=\"5=4, 1\"; leftNew=\"31\"; rightNew=\"4=7, 3\" },\n  @{ row=5; col=2; leftOld=\"61\"; rightOld=\"4=15, 1\"; leftNew=\"64\"; rightNew=\"6=10, 4\" },\n  @{ row=5; col=3; leftOld=\"54\"; rightOld=\"7=7, 5\"; leftNew=\"74\"; rightNew=\"4=18, 2\" },\n  @{ row=5; col=4; leftOld=\"40\"; rightOld=\"2=20, 0\"; leftNew=\"53\"; rightNew=\"8=6, 5\" },\n  @{ row=5; col=5; leftOld=\"79\"; rightOld=\"3=26, 1\"; leftNew=\"59\"; rightNew=\"9=6, 5\" },\n  @{ row=9; col=1; leftOld=\"93\"; rightOld=\"7=13, 2\"; leftNew=\"22\"; rightNew=\"6=3, 4\" },\n  @{ row=9; col=2; leftOld=\"78\"; rightOld=\"3=26, 0\"; leftNew=\"57\"; rightNew=\"7=8, 1\" },\n  @{ row=9; col=3; leftOld=\"97\"; rightOld=\"3=32, 1\"; leftNew=\"46\"; rightNew=\"5=9, 1\" },\n  @{ row=9; col=4; leftOld=\"69\"; rightOld=\"3=23, 0\"; leftNew=\"78\"; rightNew=\"2=39, 0\" },\n  @{ row=9; col=5; leftOld=\"88\"; rightOld=\"3=29, 1\"; leftNew=\"20\"; rightNew=\"9=2, 2\" },\n  @{ row=13; col=1; leftOld=\"17\"; rightOld=\"3=5, 2\"; leftNew=\"82\"; rightNew=\"9=9, 1\" },\n  @{ row=13; col=2; leftOld=\"52\"; rightOld=\"4=13, 0\"; leftNew=\"58\"; rightNew=\"2=29, 0\" },\n  @{ row=13; col=3; leftOld=\"49\"; rightOld=\"2=24, 1\"; leftNew=\"28\"; rightNew=\"8=3, 4\" },\n  @{ row=13; col=4; leftOld=\"31\"; rightOld=\"6=5, 1\"; leftNew=\"29\"; rightNew=\"2=14, 1\" },\n  @{ row=13; col=5; leftOld=\"61\"; rightOld=\"3=20, 1\"; leftNew=\"86\"; rightNew=\"3=28, 2\" },\n  @{ row=17; col=1; leftOld=\"29\"; rightOld=\"5=5, 4\"; leftNew=\"77\"; rightNew=\"4=19, 1\" },\n  @{ row=17; col=2; leftOld=\"89\"; rightOld=\"9=9, 8\"; leftNew=\"13\"; rightNew=\"9=1, 4\" },\n  @{ row=17; col=3; leftOld=\"89\"; rightOld=\"5=17, 4\"; leftNew=\"71\"; rightNew=\"6=11, 5\" },\n  @{ row=17; col=4; leftOld=\"93\"; rightOld=\"7=13, 2\"; leftNew=\"68\"; rightNew=\"5=13, 3\" },\n  @{ row=17; col=5; leftOld=\"98\"; rightOld=\"2=49, 0\"; leftNew=\"42\"; rightNew=\"6=7, 0\" }\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($i = 0; $i -lt $entries.Count; $i++) {\n  $entry = $entries[$i]\n  $cell = $t.Cell($entry.row, $entry.col)\n  $oldText = $entry.leftOld + $div + $entry.rightOld\n  $newText = $entry.leftNew + $div + $entry.rightNew\n  $actual = $cell.Range.Text\n  $expectedLen = $oldText.Length + 2\n  if ($actual.Length -ne $expectedLen) {\n    throw \"Cell ($($entry.row),$($entry.col)): expected length $expectedLen but found $($actual.Length)\"\n  }\n  if (-not $actual.Contains($entry.leftOld)) {\n    throw \"Cell ($($entry.row),$($entry.col)): did not contain expected left part '$($entry.leftOld)'\"\n  }\n  if (-not $actual.Contains($entry.rightOld)) {\n    throw \"Cell ($($entry.row),$($entry.col)): did not contain expected right part '$($entry.rightOld)'\"\n  }\n  $cell.Range.Text = $newText\n}\n"}
